# Applies the diff:
#  - Oval 46 (id=47): clears the "Slope " run in its 3rd paragraph, leaving an
#    empty paragraph (endParaRPr only).
#  - Oval 154 (id=155): clears the "Intercept" run in its 3rd paragraph, while
#    keeping the trailing " " run.
#  - Deletes the two "I-Mean" / "S-Mean" TextBox 319 shapes.
#  - Repositions "TextBox 9" (id=10) and "TextBox 61" (id=62).
#  - Adds four new centered textboxes: "Slope", "M = S-Mean", "M = I-Mean",
#    "Intercept" (duplicated from an existing similarly-formatted textbox so
#    they inherit the same run/paragraph/bodyPr formatting).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Find-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

function Find-ShapeByText($shapes, $text) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.HasTextFrame) {
            if ($candidate.TextFrame.HasText) {
                if ($candidate.TextFrame.TextRange.Text -eq $text) {
                    return $candidate
                }
            }
        }
    }
    return $null
}

# 1) "Slope " label on the slope Oval loses its text run (paragraph 3 of 3).
$slopeOval = Find-ShapeById $s.Shapes 47
$slopeOval.TextFrame.TextRange.Paragraphs(3).Text = ""

# 2) "Intercept" run removed from the intercept Oval, " " run kept
#    (paragraph 3 of 4, which has 2 runs: "Intercept" then " ").
$interceptOval = Find-ShapeById $s.Shapes 155
$interceptOval.TextFrame.TextRange.Paragraphs(3).Runs(1).Text = ""

# 3) Remove the stray "I-Mean" / "S-Mean" textboxes entirely.
$iMean = Find-ShapeByText $s.Shapes "I-Mean"
$iMean.Delete()
$sMean = Find-ShapeByText $s.Shapes "S-Mean"
$sMean.Delete()

# 4) Reposition the two "Construct Name" textboxes that sit beside the ovals.
$constructBottom = Find-ShapeById $s.Shapes 10
$constructBottom.Left = 59.206
$constructBottom.Top = 234.943

$constructTop = Find-ShapeById $s.Shapes 62
$constructTop.Left = 67.9386
$constructTop.Top = 65.9088

# 5) Add the four new labels, duplicating an existing centered textbox so the
#    run/paragraph/bodyPr formatting (spAutoFit, rtlCol, dirty flags, etc.)
#    matches the rest of the deck.
$template = Find-ShapeById $s.Shapes 10

$slopeLabel = $template.Duplicate().Item(1)
$slopeLabel.Name = "TextBox 62"
$slopeLabel.Left = 64.97402
$slopeLabel.Top = 85.2085
$slopeLabel.Width = 170.9915
$slopeLabel.Height = 29.0813
$slopeLabel.TextFrame.TextRange.Text = "Slope"

$sMeanLabel = $template.Duplicate().Item(1)
$sMeanLabel.Name = "TextBox 63"
$sMeanLabel.Left = 64.97402
$sMeanLabel.Top = 104.2937
$sMeanLabel.Width = 170.9915
$sMeanLabel.Height = 29.0813
$sMeanLabel.TextFrame.TextRange.Text = "M = S-Mean"

$iMeanLabel = $template.Duplicate().Item(1)
$iMeanLabel.Name = "TextBox 64"
$iMeanLabel.Left = 60.1029
$iMeanLabel.Top = 271.463
$iMeanLabel.Width = 170.9915
$iMeanLabel.Height = 29.0813
$iMeanLabel.TextFrame.TextRange.Text = "M = I-Mean"

$interceptLabel = $template.Duplicate().Item(1)
$interceptLabel.Name = "TextBox 65"
$interceptLabel.Left = 62.6703
$interceptLabel.Top = 252.49992
$interceptLabel.Width = 170.9915
$interceptLabel.Height = 29.0813
$interceptLabel.TextFrame.TextRange.Text = "Intercept"
